$wb = $excel.ActiveWorkbook

# --- Sheet "By year" ---
$ws1 = $wb.Worksheets.Item("By year")

$ws1.Range("G2").Value = 71.01
$ws1.Range("H2").Value = 504.6
$ws1.Range("I2").Value = 961.4
$ws1.Range("J2").Value = 0.1649
$ws1.Range("K2").Value = 836.1
$ws1.Range("L2").Value = 1214
$ws1.Range("O2").Value = 0.1451
$ws1.Range("P2").Value = 0.202

$ws1.Range("G3").Value = 80.86
$ws1.Range("H3").Value = 290.5
$ws1.Range("I3").Value = 758.9
$ws1.Range("J3").Value = 0.1284
$ws1.Range("K3").Value = 630.9
$ws1.Range("L3").Value = 951.3
$ws1.Range("O3").Value = 0.1111
$ws1.Range("P3").Value = 0.1552

$ws1.Range("G4").Value = 71.01
$ws1.Range("H4").Value = 310.2
$ws1.Range("I4").Value = 614.6
$ws1.Range("J4").Value = 0.1054
$ws1.Range("K4").Value = 472.5
$ws1.Range("L4").Value = 807.3
$ws1.Range("O4").Value = 0.0789
$ws1.Range("P4").Value = 0.1285

# --- Sheet "By year, livestock cat." ---
$ws2 = $wb.Worksheets.Item("By year, livestock cat.")

$ws2.Range("H2").Value = 71.01
$ws2.Range("I2").Value = 257.1
$ws2.Range("J2").Value = 328.1
$ws2.Range("K2").Value = 0.3066
$ws2.Range("L2").Value = 208.4
$ws2.Range("M2").Value = 432.9
$ws2.Range("P2").Value = 0.1942
$ws2.Range("Q2").Value = 0.3882

$ws2.Range("H3").Value = 128.7
$ws2.Range("I3").Value = 504.6
$ws2.Range("J3").Value = 633.3
$ws2.Range("K3").Value = 0.1331
$ws2.Range("L3").Value = 543.5
$ws2.Range("M3").Value = 862.1
$ws2.Range("P3").Value = 0.1199
$ws2.Range("Q3").Value = 0.1799

$ws2.Range("H4").Value = 80.86
$ws2.Range("I4").Value = 266.2
$ws2.Range("J4").Value = 347
$ws2.Range("K4").Value = 0.2444
$ws2.Range("L4").Value = 254.6
$ws2.Range("M4").Value = 466.9
$ws2.Range("P4").Value = 0.1766
$ws2.Range("Q4").Value = 0.32

$ws2.Range("H5").Value = 121.3
$ws2.Range("I5").Value = 290.5
$ws2.Range("J5").Value = 411.9
$ws2.Range("K5").Value = 0.09173
$ws2.Range("L5").Value = 317.1
$ws2.Range("M5").Value = 553.5
$ws2.Range("P5").Value = 0.0768
$ws2.Range("Q5").Value = 0.1184

$ws2.Range("H6").Value = 71.01
$ws2.Range("I6").Value = 139.3
$ws2.Range("J6").Value = 210.3
$ws2.Range("K6").Value = 0.1966
$ws2.Range("L6").Value = 162.3
$ws2.Range("M6").Value = 246.1
$ws2.Range("P6").Value = 0.149
$ws2.Range("Q6").Value = 0.2174

$ws2.Range("H7").Value = 94.03
$ws2.Range("I7").Value = 310.2
$ws2.Range("J7").Value = 404.2
$ws2.Range("K7").Value = 0.08492
$ws2.Range("L7").Value = 291
$ws2.Range("M7").Value = 573.1
$ws2.Range("P7").Value = 0.06158
$ws2.Range("Q7").Value = 0.1145

# --- Sheet "By application" ---
$ws3 = $wb.Worksheets.Item("By application")

$ws3.Range("AB2").Value = 18.26
$ws3.Range("AC2").Value = 0.2029

$ws3.Range("AB3").Value = 32.14
$ws3.Range("AC3").Value = 0.3571

$ws3.Range("AB4").Value = 12.46
$ws3.Range("AC4").Value = 0.1038

$ws3.Range("AB5").Value = 17.2
$ws3.Range("AC5").Value = 0.1434

$ws3.Range("AB6").Value = 20.22
$ws3.Range("AC6").Value = 0.1925

$ws3.Range("AB7").Value = 27.95
$ws3.Range("AC7").Value = 0.2662

$ws3.Range("AB8").Value = 13.77
$ws3.Range("AC8").Value = 0.102

$ws3.Range("AB9").Value = 10.56
$ws3.Range("AC9").Value = 0.08804

$ws3.Range("AB10").Value = 18.26
$ws3.Range("AC10").Value = 0.2029

$ws3.Range("AB11").Value = 17.42
$ws3.Range("AC11").Value = 0.1935

$ws3.Range("AB12").Value = 9.1
$ws3.Range("AC12").Value = 0.07583

$ws3.Range("AB13").Value = 10.58
$ws3.Range("AC13").Value = 0.08813
